# Fixed Transformer and List Initializations
#
# - Transformers sheet gains an "Input" column (right after Name) and a
#   trailing "etc" column, and its existing headers shuffle accordingly
#   (Capex/Opex move up next to Input, the efficiency columns follow).
# - The previously-active tab (Connectors) is no longer the selected one;
#   Transformers becomes the active/selected sheet instead.

$wb = $excel.ActiveWorkbook

$wsTransformers = $wb.Worksheets.Item("Transformers")

# Rebuild the header row (row 1) of the Transformers sheet with the new
# column layout: Name, Input, Capex, Opex, TotalEff, Prod1, SubEff1,
# Prod2, SubEff2, Prod3, SubEff3, etc.
$wsTransformers.Cells.Item(1, 1).Value = "Name"
$wsTransformers.Cells.Item(1, 2).Value = "Input"
$wsTransformers.Cells.Item(1, 3).Value = " Capex"
$wsTransformers.Cells.Item(1, 4).Value = "Opex"
$wsTransformers.Cells.Item(1, 5).Value = "TotalEff"
$wsTransformers.Cells.Item(1, 6).Value = "Prod1"
$wsTransformers.Cells.Item(1, 7).Value = "SubEff1"
$wsTransformers.Cells.Item(1, 8).Value = "Prod2"
$wsTransformers.Cells.Item(1, 9).Value = "SubEff2"
$wsTransformers.Cells.Item(1, 10).Value = "Prod3"
$wsTransformers.Cells.Item(1, 11).Value = "SubEff3"
$wsTransformers.Cells.Item(1, 12).Value = "etc"

# Make Transformers the active sheet/tab, with L7 as the selected cell,
# which also clears tabSelected from whichever sheet (Connectors) used to
# be active.
$wsTransformers.Activate()
$wsTransformers.Range("L7").Select()
